{"js": "// Update the title/date paragraph (first paragraph in the document body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2023-04-23 Sunday\", Word.InsertLocation.replace);\n\n// Update every multiplication-table cell (20 rows x 5 cols) in row-major order\n// by overwriting the table's `values` grid, which preserves each cell's existing\n// run formatting (font/size) while only changing the text.\nconst newValues = [\n  [\"77\u00d730=2310\", \"93\u00d795=8835\", \"73\u00d735=2555\", \"65\u00d760=3900\", \"31\u00d719=589\"],\n  [\"15\u00d746=690\", \"70\u00d762=4340\", \"18\u00d740=720\", \"49\u00d715=735\", \"45\u00d727=1215\"],\n  [\"41\u00d739=1599\", \"98\u00d766=6468\", \"24\u00d738=912\", \"97\u00d716=1552\", \"92\u00d736=3312\"],\n  [\"42\u00d792=3864\", \"31\u00d732=992\", \"99\u00d723=2277\", \"52\u00d746=2392\", \"23\u00d775=1725\"],\n  [\"53\u00d780=4240\", \"24\u00d741=984\", \"26\u00d777=2002\", \"13\u00d743=559\", \"77\u00d739=3003\"],\n  [\"81\u00d790=7290\", \"49\u00d765=3185\", \"86\u00d790=7740\", \"33\u00d766=2178\", \"25\u00d726=650\"],\n  [\"32\u00d752=1664\", \"68\u00d772=4896\", \"91\u00d774=6734\", \"65\u00d746=2990\", \"33\u00d710=330\"],\n  [\"50\u00d774=3700\", \"79\u00d746=3634\", \"55\u00d741=2255\", \"67\u00d795=6365\", \"74\u00d785=6290\"],\n  [\"98\u00d749=4802\", \"40\u00d782=3280\", \"92\u00d783=7636\", \"23\u00d733=759\", \"65\u00d776=4940\"],\n  [\"42\u00d763=2646\", \"36\u00d768=2448\", \"84\u00d777=6468\", \"86\u00d7100=8600\", \"20\u00d797=1940\"],\n  [\"18\u00d797=1746\", \"27\u00d766=1782\", \"52\u00d753=2756\", \"31\u00d769=2139\", \"71\u00d719=1349\"],\n  [\"91\u00d789=8099\", \"10\u00d764=640\", \"63\u00d745=2835\", \"32\u00d751=1632\", \"85\u00d737=3145\"],\n  [\"63\u00d783=5229\", \"56\u00d711=616\", \"37\u00d757=2109\", \"21\u00d795=1995\", \"10\u00d770=700\"],\n  [\"47\u00d762=2914\", \"95\u00d793=8835\", \"13\u00d785=1105\", \"90\u00d730=2700\", \"88\u00d797=8536\"],\n  [\"46\u00d767=3082\", \"66\u00d747=3102\", \"25\u00d718=450\", \"80\u00d714=1120\", \"34\u00d748=1632\"],\n  [\"17\u00d738=646\", \"38\u00d741=1558\", \"57\u00d726=1482\", \"68\u00d781=5508\", \"25\u00d741=1025\"],\n  [\"93\u00d7100=9300\", \"64\u00d719=1216\", \"55\u00d757=3135\", \"37\u00d790=3330\", \"44\u00d719=836\"],\n  [\"74\u00d753=3922\", \"99\u00d738=3762\", \"11\u00d775=825\", \"38\u00d729=1102\", \"87\u00d733=2871\"],\n  [\"87\u00d774=6438\", \"82\u00d717=1394\", \"13\u00d723=299\", \"85\u00d746=3910\", \"43\u00d714=602\"],\n  [\"68\u00d766=4488\", \"40\u00d735=1400\", \"47\u00d753=2491\", \"41\u00d774=3034\", \"93\u00d7100=9300\"]\n];\n\nconst table = context.document.body.tables.getFirst();\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the title/date paragraph (first paragraph in the document body).\n$d.Paragraphs.Item(1).Range.Text = \"2023-04-23 Sunday\"\n\n# Update each multiplication-table cell in row-major order (row 1..20, col 1..5).\n$newValues = @(\n    \"77\u00d730=2310\",\n    \"93\u00d795=8835\",\n    \"73\u00d735=2555\",\n    \"65\u00d760=3900\",\n    \"31\u00d719=589\",\n    \"15\u00d746=690\",\n    \"70\u00d762=4340\",\n    \"18\u00d740=720\",\n    \"49\u00d715=735\",\n    \"45\u00d727=1215\",\n    \"41\u00d739=1599\",\n    \"98\u00d766=6468\",\n    \"24\u00d738=912\",\n    \"97\u00d716=1552\",\n    \"92\u00d736=3312\",\n    \"42\u00d792=3864\",\n    \"31\u00d732=992\",\n    \"99\u00d723=2277\",\n    \"52\u00d746=2392\",\n    \"23\u00d775=1725\",\n    \"53\u00d780=4240\",\n    \"24\u00d741=984\",\n    \"26\u00d777=2002\",\n    \"13\u00d743=559\",\n    \"77\u00d739=3003\",\n    \"81\u00d790=7290\",\n    \"49\u00d765=3185\",\n    \"86\u00d790=7740\",\n    \"33\u00d766=2178\",\n    \"25\u00d726=650\",\n    \"32\u00d752=1664\",\n    \"68\u00d772=4896\",\n    \"91\u00d774=6734\",\n    \"65\u00d746=2990\",\n    \"33\u00d710=330\",\n    \"50\u00d774=3700\",\n    \"79\u00d746=3634\",\n    \"55\u00d741=2255\",\n    \"67\u00d795=6365\",\n    \"74\u00d785=6290\",\n    \"98\u00d749=4802\",\n    \"40\u00d782=3280\",\n    \"92\u00d783=7636\",\n    \"23\u00d733=759\",\n    \"65\u00d776=4940\",\n    \"42\u00d763=2646\",\n    \"36\u00d768=2448\",\n    \"84\u00d777=6468\",\n    \"86\u00d7100=8600\",\n    \"20\u00d797=1940\",\n    \"18\u00d797=1746\",\n    \"27\u00d766=1782\",\n    \"52\u00d753=2756\",\n    \"31\u00d769=2139\",\n    \"71\u00d719=1349\",\n    \"91\u00d789=8099\",\n    \"10\u00d764=640\",\n    \"63\u00d745=2835\",\n    \"32\u00d751=1632\",\n    \"85\u00d737=3145\",\n    \"63\u00d783=5229\",\n    \"56\u00d711=616\",\n    \"37\u00d757=2109\",\n    \"21\u00d795=1995\",\n    \"10\u00d770=700\",\n    \"47\u00d762=2914\",\n    \"95\u00d793=8835\",\n    \"13\u00d785=1105\",\n    \"90\u00d730=2700\",\n    \"88\u00d797=8536\",\n    \"46\u00d767=3082\",\n    \"66\u00d747=3102\",\n    \"25\u00d718=450\",\n    \"80\u00d714=1120\",\n    \"34\u00d748=1632\",\n    \"17\u00d738=646\",\n    \"38\u00d741=1558\",\n    \"57\u00d726=1482\",\n    \"68\u00d781=5508\",\n    \"25\u00d741=1025\",\n    \"93\u00d7100=9300\",\n    \"64\u00d719=1216\",\n    \"55\u00d757=3135\",\n    \"37\u00d790=3330\",\n    \"44\u00d719=836\",\n    \"74\u00d753=3922\",\n    \"99\u00d738=3762\",\n    \"11\u00d775=825\",\n    \"38\u00d729=1102\",\n    \"87\u00d733=2871\",\n    \"87\u00d774=6438\",\n    \"82\u00d717=1394\",\n    \"13\u00d723=299\",\n    \"85\u00d746=3910\",\n    \"43\u00d714=602\",\n    \"68\u00d766=4488\",\n    \"40\u00d735=1400\",\n    \"47\u00d753=2491\",\n    \"41\u00d774=3034\",\n    \"93\u00d7100=9300\"\n)\n\n$table = $d.Tables.Item(1)\n$cols = 5\n$idx = 0\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $table.Cell($r, $c).Range.Text = $newValues[$idx]\n        $idx = $idx + 1\n    }\n}\n\nWrite-Output \"Updated $idx table cells.\"\n"}
